$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.902300119400024
$ws.Range("B1").Value = 2.253734350204468
$ws.Range("C1").Value = 2.436330556869507
$ws.Range("D1").Value = 3.459383487701416
$ws.Range("E1").Value = 1.22719669342041
